# Update the cryptos list with the latest scraped prices / 1h volume changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "30.647.15"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value2 = "  +0.60%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.881.51"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value2 = "  -0.48%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value2 = "  +0.20%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "239.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value2 = "  -0.10%  "

# Row 6 - USDC
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "1.000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value2 = "  +0.11%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.4830"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value2 = "  -0.31%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "1.882.78"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value2 = "  -0.22%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.2839"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value2 = "  -1.92%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.06540"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value2 = "  -1.33%  "

# Row 11 - WrappedEther
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "1.904.96"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value2 = "  +0.87%  "

# Row 12 - TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.07514"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value2 = "  +1.23%  "

# Row 13 - Solana
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "16.63"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value2 = "  -2.22%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "5.111"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value2 = "  -1.54%  "

# Row 15 - Litecoin
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "88.95"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value2 = "  -0.19%  "

# Row 16 - Polygon
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.6668"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value2 = "  +0.22%  "

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "30.602.14"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value2 = "  +0.66%  "

# Row 18 - WrappedliquidstakedEther2.0
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "2.285.34"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value2 = "  +7.23%  "

# Row 19 - Avalanche
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "13.37"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value2 = "  -1.28%  "

# Row 20 - Dai
$ws.Range("E20").Value2 = "  +0.08%  "

# Row 21 - ShibaInu
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "0.000007624"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value2 = "  -2.02%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "232.06"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value2 = "  +8.64%  "

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "5.301"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value2 = "  -2.56%  "

# Row 24 - BinanceUSD
$ws.Range("E24").Value2 = "  +0.20%  "

# Row 25 - Chainlink
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "6.190"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value2 = "  -0.26%  "

# Row 26 - Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "9.371"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value2 = "  -0.80%  "

# Row 27 - Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "167.72"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value2 = "  +1.42%  "

# Row 28 - EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "18.67"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value2 = "  +0.49%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "1.949"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value2 = "  +0.15%  "

# Row 30 - Toncoin
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.421"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value2 = "  -1.12%  "

# Row 31 - Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.09566"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value2 = "  +4.05%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "4.377"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value2 = "  +1.32%  "

# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "4.050"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value2 = "  -0.87%  "

# Row 34 - Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.05041"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value2 = "  -0.79%  "

# Row 35 - ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.215"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value2 = "  +4.59%  "

# Row 36 - ImmutableX
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.7499"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value2 = "  -0.51%  "

# Row 37 - HuobiToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "2.699"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value2 = "  -0.17%  "

# Row 38 - VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.01856"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value2 = "  -2.45%  "

# Row 39 - MXToken
$ws.Range("E39").Value2 = "  -0.58%  "

# Row 40 - RenderToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "2.104"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value2 = "  +0.46%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.9177"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value2 = "  -0.20%  "

# Row 42 - Quant
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "106.72"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value2 = "  -0.65%  "

# Row 43 - TheSandbox
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.4300"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value2 = "  -1.30%  "

# Row 44 - FraxShare
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "5.829"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value2 = "  -4.32%  "

# Row 45 - PaxDollar
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "1.001"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value2 = "  -0.07%  "

# Row 46 - Aptos
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "7.436"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value2 = "  -3.01%  "

# Row 49 - NEARProtocol
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "1.483"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value2 = "  -6.66%  "

# Row 50 - EnergySwap
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "8.933"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value2 = "  -0.21%  "

# Row 51 - Elrond
$ws.Range("E51").Value2 = "  -1.56%  "

# Rows 47 and 48 swap places: Algorand now ranked 47 (was Aave), Aave now ranked 48 (was Algorand)
# Row 47 - Algorand (was Aave)
$ws.Range("B47").Value2 = "Algorand"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.1288"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value2 = "  -4.61%  "

# Row 48 - Aave (was Algorand)
$ws.Range("B48").Value2 = "Aave"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "64.35"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value2 = "  -2.50%  "

